$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.200.39'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').Value = '2.633.25'
$ws.Range('E3').Value = '  +4.28%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '523.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.13%  '
$ws.Range('E7').Value = '  -0.45%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.571'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.17%  '
$ws.Range('D9').Value = '2.658.24'
$ws.Range('E9').Value = '  +4.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.33'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.106'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.97%  '
$ws.Range('E12').Value = '  +3.32%  '
$ws.Range('E13').Value = '  -1.12%  '
$ws.Range('D14').Value = '3.099.83'
$ws.Range('E14').Value = '  +4.14%  '
$ws.Range('D15').Value = '59.129.12'
$ws.Range('E15').Value = '  +2.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.06'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('E17').Value = '  +3.45%  '
$ws.Range('D18').Value = '2.652.32'
$ws.Range('E18').Value = '  +4.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '348.68'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.52%  '
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('E21').Value = '  +3.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.18'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.09%  '
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.80'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.70%  '
$ws.Range('E25').Value = '  +3.40%  '
$ws.Range('E26').Value = '  +4.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.993'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.57%  '
$ws.Range('D28').Value = '0.0₃0812'
$ws.Range('E28').Value = '  +5.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.14'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.998'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.28'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.77%  '
$ws.Range('E32').Value = '  +4.80%  '
$ws.Range('E33').Value = '  +3.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '149.97'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.979'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.86%  '
$ws.Range('E36').Value = '  +4.65%  '
$ws.Range('E37').Value = '  +3.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '36.82'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.78%  '
$ws.Range('E39').Value = '  +5.69%  '
$ws.Range('E40').Value = '  +6.17%  '
$ws.Range('E41').Value = '  +4.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '279.80'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.31%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.611'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.55%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0988'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.35%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.994'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.64'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.64%  '
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.76'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0230'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.30'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('D51').Value = '1.993.50'
$ws.Range('E51').Value = '  +5.75%  '
